$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated NATMI Col3a1-Ddr2 edge statistics (Natmi following Dr Hou advice).
# Ligand-expressing cells (E) and Receptor-expressing cells (K) moved from 1 to 3,
# which changes every downstream derived expression/specificity value in columns
# E, G-K, M-T for data rows 2-17.

$rowData = @{
    2 = @{ "E" = 3.0; "G" = 21.62966433333333; "H" = 64.888993; "I" = 0.004276908378962984; "J" = 0.004276908378962984; "K" = 3.0; "M" = 6.546185; "N" = 19.638555; "O" = 0.06829436374455893; "P" = 0.06829436374455893; "Q" = 141.5917842139017; "R" = 1274.326057925115; "S" = 0.0002920887365350499; "T" = 0.0002920887365350499 }
    3 = @{ "E" = 3.0; "G" = 21.62966433333333; "H" = 64.888993; "I" = 0.004276908378962984; "J" = 0.004276908378962984; "K" = 3.0; "M" = 71.03440333333334; "N" = 213.10321; "O" = 0.7410803971510699; "P" = 0.7410803971510698; "Q" = 1536.450300218615; "R" = 13828.05270196753; "S" = 0.003169532960060627; "T" = 0.003169532960060626 }
    4 = @{ "E" = 3.0; "G" = 21.62966433333333; "H" = 64.888993; "I" = 0.004276908378962984; "J" = 0.004276908378962984; "K" = 3.0; "M" = 0.3068453333333334; "N" = 0.920536; "O" = 0.003201224347919758; "P" = 0.003201224347919758; "Q" = 6.636961562249779; "R" = 59.732654060248; "S" = 0.00001369134323655833; "T" = 0.00001369134323655832 }
    5 = @{ "E" = 3.0; "G" = 21.62966433333333; "H" = 64.888993; "I" = 0.004276908378962984; "J" = 0.004276908378962984; "K" = 3.0; "M" = 17.965059; "N" = 53.895177; "O" = 0.1874240147564516; "P" = 0.1874240147564516; "Q" = 388.578195898529; "R" = 3497.203763086761; "S" = 0.0008015953391307497; "T" = 0.0008015953391307497 }
    6 = @{ "E" = 3.0; "G" = 4967.017741; "H" = 14901.053223; "I" = 0.9821456064948035; "J" = 0.9821456064948036; "K" = 3.0; "M" = 6.546185; "N" = 19.638555; "O" = 0.06829436374455893; "P" = 0.06829436374455893; "Q" = 32515.01703086809; "R" = 292635.1532778128; "S" = 0.06707500930007655; "T" = 0.06707500930007655 }
    7 = @{ "E" = 3.0; "G" = 4967.017741; "H" = 14901.053223; "I" = 0.9821456064948035; "J" = 0.9821456064948036; "K" = 3.0; "M" = 71.03440333333334; "N" = 213.10321; "O" = 0.7410803971510699; "P" = 0.7410803971510698; "Q" = 352829.1415780162; "R" = 3175462.274202146; "S" = 0.7278488561213474; "T" = 0.7278488561213474 }
    8 = @{ "E" = 3.0; "G" = 4967.017741; "H" = 14901.053223; "I" = 0.9821456064948035; "J" = 0.9821456064948036; "K" = 3.0; "M" = 0.3068453333333334; "N" = 0.920536; "O" = 0.003201224347919758; "P" = 0.003201224347919758; "Q" = 1524.106214409725; "R" = 13716.95592968753; "S" = 0.003144068428713583; "T" = 0.003144068428713583 }
    9 = @{ "E" = 3.0; "G" = 4967.017741; "H" = 14901.053223; "I" = 0.9821456064948035; "J" = 0.9821456064948036; "K" = 3.0; "M" = 17.965059; "N" = 53.895177; "O" = 0.1874240147564516; "P" = 0.1874240147564516; "Q" = 89232.76677111171; "R" = 803094.9009400054; "S" = 0.1840776726446662; "T" = 0.1840776726446662 }
    10 = @{ "E" = 3.0; "G" = 2.249417666666667; "H" = 6.748253; "I" = 0.0004447851394313067; "J" = 0.0004447851394313068; "K" = 3.0; "M" = 6.546185; "N" = 19.638555; "O" = 0.06829436374455893; "P" = 0.06829436374455893; "Q" = 14.72510418826833; "R" = 132.525937694415; "S" = 0.00003037631810049602; "T" = 0.00003037631810049603 }
    11 = @{ "E" = 3.0; "G" = 2.249417666666667; "H" = 6.748253; "I" = 0.0004447851394313067; "J" = 0.0004447851394313068; "K" = 3.0; "M" = 71.03440333333334; "N" = 213.10321; "O" = 0.7410803971510699; "P" = 0.7410803971510698; "Q" = 159.7860417991256; "R" = 1438.07437619213; "S" = 0.0003296215477766468; "T" = 0.0003296215477766468 }
    12 = @{ "E" = 3.0; "G" = 2.249417666666667; "H" = 6.748253; "I" = 0.0004447851394313067; "J" = 0.0004447851394313068; "K" = 3.0; "M" = 0.3068453333333334; "N" = 0.920536; "O" = 0.003201224347919758; "P" = 0.003201224347919758; "Q" = 0.6902233137342222; "R" = 6.212009823608001; "S" = 0.000001423857017940383; "T" = 0.000001423857017940383 }
    13 = @{ "E" = 3.0; "G" = 2.249417666666667; "H" = 6.748253; "I" = 0.0004447851394313067; "J" = 0.0004447851394313068; "K" = 3.0; "M" = 17.965059; "N" = 53.895177; "O" = 0.1874240147564516; "P" = 0.1874240147564516; "Q" = 40.410921097309; "R" = 363.698289875781; "S" = 0.00008336341653622361; "T" = 0.00008336341653622363 }
    14 = @{ "E" = 3.0; "G" = 66.41617433333333; "H" = 199.248523; "I" = 0.01313269998680205; "J" = 0.01313269998680205; "K" = 3.0; "M" = 6.546185; "N" = 19.638555; "O" = 0.06829436374455893; "P" = 0.06829436374455893; "Q" = 434.7725641782517; "R" = 3912.953077604265; "S" = 0.0008968893898468237; "T" = 0.0008968893898468237 }
    15 = @{ "E" = 3.0; "G" = 66.41617433333333; "H" = 199.248523; "I" = 0.01313269998680205; "J" = 0.01313269998680205; "K" = 3.0; "M" = 71.03440333333334; "N" = 213.10321; "O" = 0.7410803971510699; "P" = 0.7410803971510698; "Q" = 4717.833315450982; "R" = 42460.49983905883; "S" = 0.009732386521885116; "T" = 0.009732386521885116 }
    16 = @{ "E" = 3.0; "G" = 66.41617433333333; "H" = 199.248523; "I" = 0.01313269998680205; "J" = 0.01313269998680205; "K" = 3.0; "M" = 0.3068453333333334; "N" = 0.920536; "O" = 0.003201224347919758; "P" = 0.003201224347919758; "Q" = 20.37949315203645; "R" = 183.415438368328; "S" = 0.00004204071895167622; "T" = 0.00004204071895167622 }
    17 = @{ "E" = 3.0; "G" = 66.41617433333333; "H" = 199.248523; "I" = 0.01313269998680205; "J" = 0.01313269998680205; "K" = 3.0; "M" = 17.965059; "N" = 53.895177; "O" = 0.1874240147564516; "P" = 0.1874240147564516; "Q" = 1193.170490452619; "R" = 10738.53441407357; "S" = 0.00246138335611844; "T" = 0.00246138335611844 }
}

foreach ($r in $rowData.Keys) {
    foreach ($col in $rowData[$r].Keys) {
        $ws.Range("$col$r").Value = $rowData[$r][$col]
    }
}
